# Insert a new weekly price record for "Macroferia Regional de Talca - Apio"
# at row 338 (pushing the existing rows 338:373 down to 339:374), then
# populate the new row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 338; existing rows shift down by one.
$ws.Rows.Item(338).Insert()

# Fill the new row 338 with the new record's data.
$ws.Cells.Item(338, 1).Value = 5
$ws.Cells.Item(338, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(338, 3).Value = "Maule"
$ws.Cells.Item(338, 4).Value = 45194
$ws.Cells.Item(338, 5).Value = 7
$ws.Cells.Item(338, 6).Value = 100112017
$ws.Cells.Item(338, 7).Value = "Apio"
$ws.Cells.Item(338, 8).Value = "Americana (o)"
$ws.Cells.Item(338, 9).Value = "Primera"
$ws.Cells.Item(338, 10).Value = 700
$ws.Cells.Item(338, 11).Value = 6500
$ws.Cells.Item(338, 12).Value = 6500
$ws.Cells.Item(338, 13).Value = 6500
$ws.Cells.Item(338, 14).Value = "`$/docena de matas"
$ws.Cells.Item(338, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(338, 16).Value = 1083
$ws.Cells.Item(338, 17).Value = 6
$ws.Cells.Item(338, 18).Value = "Hortaliza"

Write-Output "Row inserted and populated."
